$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A, shifting B:H -> C:I
$ws.Columns("A").Insert()

# New header / data for the inserted column A
$ws.Range("A1").Value = "Sopimuksen SAP-tunnus"
$ws.Range("A2").Value = "fakeSapContractId"

# The "contact" SAP id was missing - fill it into column B (where the old
# farmer id used to live before the column insert shifted it here), with
# its own explicit text-formatted style
$ws.Range("B2").Value = "fakeSapContactId"
$ws.Range("B2").Font.ColorIndex = 1
$ws.Range("B2").NumberFormat = "@"

# Delivery place SAP id placeholder (was "01")
$ws.Range("G2").Value = "fakeDPSapId"

# Restore column A's intended width (engine snaps to a ~1/6-char grid, so
# feed it a value that lands on the closest achievable width to 24.45)
$ws.Columns("A").ColumnWidth = 23.6
